# Apply trade #2 close update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1299.96
$summary.Range("B4").Value = -0.04
$summary.Range("B5").Value = -0.4
$summary.Range("B6").Value = 2
$summary.Range("B8").Value = 2

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95999999999999
$status.Range("D4").Value = 2
$status.Range("E4").Value = -0.04
$status.Range("F4").Value = -0.04

# --- Helper to append the new trade row to a trades-log sheet ---
function Add-TradeRow($sheet) {
    $sheet.Range("A3").Value = 2

    # Force the date-looking string to stay as literal text instead of
    # being auto-converted to a date serial number, then restore the
    # default cell style so no stray formatting is left behind.
    $sheet.Range("B3").NumberFormat = "@"
    $sheet.Range("B3").Value = "2026-02-17"
    $sheet.Range("B3").Style = "Normal"

    $sheet.Range("C3").Value = "19:43:37"
    $sheet.Range("D3").Value = "MarketMaking"
    $sheet.Range("E3").Value = "UP"
    $sheet.Range("F3").Value = 0.4
    $sheet.Range("G3").Value = 0.37
    $sheet.Range("H3").Value = "CLOSED"
    $sheet.Range("I3").Value = -7.5
    $sheet.Range("J3").Value = -0.03
    $sheet.Range("K3").Value = 99.95999999999999
    $sheet.Range("L3").Value = 0
    $sheet.Range("M3").Value = 0
    $sheet.Range("N3").Value = 0.6
    $sheet.Range("O3").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P3").Value = "early_exit"
    $sheet.Range("Q3").Value = 0.13
}

# --- All Trades sheet ---
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# --- MarketMaking sheet ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
